# Remove the trailing site-chrome paragraphs that were scraped along with
# the bibliography entry ("Ver no Jupiter Salvar em pdf Salvar em docx" and
# the "© 2020 ..." footer line), plus the now-superfluous blank paragraph
# that separated them from the bibliography text.
$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $i
    }
    if ($text -like "*Powered by Jekyll and Github pages*") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    # Also swallow the blank paragraph immediately preceding the
    # "Ver no Jupiter..." paragraph, so only a single blank paragraph is
    # left before the trailing page-break paragraph.
    $deleteFrom = $d.Paragraphs.Item($startPara - 1).Range.Start
    $deleteTo = $d.Paragraphs.Item($endPara).Range.End

    $killRange = $d.Range($deleteFrom, $deleteTo)
    $killRange.Delete()
}
